$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "Abrade / Hour of Devastation" line (old row 19) -- everything
# below it shifts up by one row.
$ws.Rows.Item(19).Delete()

# Insert two new rows just above "Sweltering Suns" (now row 24) to record the
# two Negate purchases, pushing that row (and the totals row) down.
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(24).Insert()

# Row 24: Negate / Magic 2013
$ws.Cells.Item(24, 1).Value = "Negate"
$ws.Cells.Item(24, 2).Value = "Magic 2013"
$ws.Cells.Item(24, 3).Value = "Normal"
$ws.Cells.Item(24, 4).Value = 0.08
$ws.Cells.Item(24, 5).Value = 1
$ws.Cells.Item(24, 6).Formula = "=D24*E24"

# Row 25: Negate / Magic 2014
$ws.Cells.Item(25, 1).Value = "Negate"
$ws.Cells.Item(25, 2).Value = "Magic 2014"
$ws.Cells.Item(25, 3).Value = "Normal"
$ws.Cells.Item(25, 4).Value = 0.1
$ws.Cells.Item(25, 5).Value = 1
$ws.Cells.Item(25, 6).Formula = "=D25*E25"

$ws.Range("F25").Select()
